$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Straightforward text/link/percentage cell updates
$ws.Range("D2").Value = "28.506.52"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.870.81"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -1.68%  "
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.871.96"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "28.529.71"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").Value = "2.087.88"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  +0.03%  "

# Price cells whose new text looks numeric: force text format so Excel
# keeps them as strings (matching the source data), then restore the
# original (unstyled) cell style so no stray number-format is left behind.
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = $style
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.47"
$ws.Range("D5").Style = $style
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = $style
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3902"
$ws.Range("D8").Style = $style
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08360"
$ws.Range("D9").Style = $style
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.105"
$ws.Range("D10").Style = $style
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.42"
$ws.Range("D13").Style = $style
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.278"
$ws.Range("D14").Style = $style
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.008"
$ws.Range("D15").Style = $style
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001103"
$ws.Range("D16").Style = $style
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.10"
$ws.Range("D17").Style = $style
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06731"
$ws.Range("D18").Style = $style
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.75"
$ws.Range("D19").Style = $style
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.917"
$ws.Range("D21").Style = $style
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.207"
$ws.Range("D24").Style = $style
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.34"
$ws.Range("D26").Style = $style
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.417"
$ws.Range("D28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "126.34"
$ws.Range("D29").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.741"
$ws.Range("D32").Style = $style
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.620"
$ws.Range("D33").Style = $style
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02458"
$ws.Range("D34").Style = $style
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06577"
$ws.Range("D35").Style = $style
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2164"
$ws.Range("D36").Style = $style
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.898"
$ws.Range("D37").Style = $style
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.027"
$ws.Range("D38").Style = $style
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.239"
$ws.Range("D39").Style = $style
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.179"
$ws.Range("D40").Style = $style
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6366"
$ws.Range("D41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("D42").Style = $style
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6007"
$ws.Range("D44").Style = $style
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.05"
$ws.Range("D45").Style = $style
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.685"
$ws.Range("D46").Style = $style
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.002"
$ws.Range("D47").Style = $style
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.37"
$ws.Range("D51").Style = $style
